# Updates LR-pair edge-weight/specificity metrics (TPM re-run) for Angptl2-Itgb1.
# Static numeric values only (no formulas in source workbook) -- apply literal replacements.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 3.440334666666667
    "H2" = 10.321004
    "I2" = 0.03090290794544385
    "J2" = 0.03090290794544386
    "M2" = 145.7007446666667
    "N2" = 437.1022340000001
    "O2" = 0.2865937750105843
    "P2" = 0.2865937750105843
    "Q2" = 501.2593228358818
    "R2" = 4511.333905522937
    "S2" = 0.008856581046889331
    "T2" = 0.008856581046889333
    "G3" = 3.440334666666667
    "H3" = 10.321004
    "I3" = 0.03090290794544385
    "J3" = 0.03090290794544386
    "O3" = 0.3320294904365841
    "P3" = 0.3320294904365841
    "Q3" = 580.7274688071592
    "R3" = 5226.547219264433
    "S3" = 0.01026067677813439
    "T3" = 0.01026067677813439
    "G4" = 3.440334666666667
    "H4" = 10.321004
    "I4" = 0.03090290794544385
    "J4" = 0.03090290794544386
    "M4" = 128.1261546666667
    "N4" = 384.378464
    "O4" = 0.2520245069956105
    "P4" = 0.2520245069956105
    "Q4" = 440.7968516064285
    "R4" = 3967.171664457856
    "S4" = 0.00778829013968122
    "T4" = 0.007788290139681221
    "G5" = 3.440334666666667
    "H5" = 10.321004
    "I5" = 0.03090290794544385
    "J5" = 0.03090290794544386
    "M5" = 65.761079
    "N5" = 197.283237
    "O5" = 0.1293522275572212
    "P5" = 0.1293522275572212
    "Q5" = 226.2401198011053
    "R5" = 2036.161078209948
    "S5" = 0.003997359980738913
    "T5" = 0.003997359980738913
    "I6" = 0.73221566931385
    "J6" = 0.7322156693138502
    "M6" = 145.7007446666667
    "N6" = 437.1022340000001
    "O6" = 0.2865937750105843
    "P6" = 0.2865937750105843
    "Q6" = 11876.87356859875
    "R6" = 106891.8621173888
    "S6" = 0.2098484527905579
    "T6" = 0.2098484527905579
    "I7" = 0.73221566931385
    "J7" = 0.7322156693138502
    "O7" = 0.3320294904365841
    "P7" = 0.3320294904365841
    "S7" = 0.24311719557196
    "T7" = 0.2431171955719601
    "I8" = 0.73221566931385
    "J8" = 0.7322156693138502
    "M8" = 128.1261546666667
    "N8" = 384.378464
    "O8" = 0.2520245069956105
    "P8" = 0.2520245069956105
    "Q8" = 10444.27153264146
    "R8" = 93998.44379377316
    "S8" = 0.184536293073284
    "T8" = 0.184536293073284
    "I9" = 0.73221566931385
    "J9" = 0.7322156693138502
    "M9" = 65.761079
    "N9" = 197.283237
    "O9" = 0.1293522275572212
    "P9" = 0.1293522275572212
    "Q9" = 5360.549273817948
    "R9" = 48244.94346436154
    "S9" = 0.09471372787804817
    "T9" = 0.0947137278780482
    "G10" = 25.13705366666667
    "H10" = 75.41116100000001
    "I10" = 0.2257943283853049
    "J10" = 0.225794328385305
    "M10" = 145.7007446666667
    "N10" = 437.1022340000001
    "O10" = 0.2865937750105843
    "P10" = 0.2865937750105843
    "Q10" = 3662.487437959298
    "R10" = 32962.38694163368
    "S10" = 0.06471124894792406
    "T10" = 0.06471124894792407
    "G11" = 25.13705366666667
    "H11" = 75.41116100000001
    "I11" = 0.2257943283853049
    "J11" = 0.225794328385305
    "O11" = 0.3320294904365841
    "P11" = 0.3320294904365841
    "Q11" = 4243.127184849378
    "R11" = 38188.14466364439
    "S11" = 0.07497037579724354
    "T11" = 0.07497037579724355
    "G12" = 25.13705366666667
    "H12" = 75.41116100000001
    "I12" = 0.2257943283853049
    "J12" = 0.225794328385305
    "M12" = 128.1261546666667
    "N12" = 384.378464
    "O12" = 0.2520245069956105
    "P12" = 0.2520245069956105
    "Q12" = 3220.714025959634
    "R12" = 28986.42623363671
    "S12" = 0.05690570429371145
    "T12" = 0.05690570429371146
    "G13" = 25.13705366666667
    "H13" = 75.41116100000001
    "I13" = 0.2257943283853049
    "J13" = 0.225794328385305
    "M13" = 65.761079
    "N13" = 197.283237
    "O13" = 0.1293522275572212
    "P13" = 0.1293522275572212
    "Q13" = 1653.039772000906
    "R13" = 14877.35794800816
    "S13" = 0.0292069993464259
    "T13" = 0.02920699934642591
    "G14" = 1.234295333333334
    "H14" = 3.702886
    "I14" = 0.01108709435540116
    "J14" = 0.01108709435540116
    "M14" = 145.7007446666667
    "N14" = 437.1022340000001
    "O14" = 0.2865937750105843
    "P14" = 0.2865937750105843
    "Q14" = 179.8377492052583
    "R14" = 1618.539742847324
    "S14" = 0.003177492225212959
    "T14" = 0.00317749222521296
    "G15" = 1.234295333333334
    "H15" = 3.702886
    "I15" = 0.01108709435540116
    "J15" = 0.01108709435540116
    "O15" = 0.3320294904365841
    "P15" = 0.3320294904365841
    "Q15" = 208.348685269521
    "R15" = 1875.138167425688
    "S15" = 0.003681242289246176
    "T15" = 0.003681242289246177
    "G16" = 1.234295333333334
    "H16" = 3.702886
    "I16" = 0.01108709435540116
    "J16" = 0.01108709435540116
    "M16" = 128.1261546666667
    "N16" = 384.378464
    "O16" = 0.2520245069956105
    "P16" = 0.2520245069956105
    "Q16" = 158.1455147830116
    "R16" = 1423.309633047104
    "S16" = 0.002794219488933793
    "T16" = 0.002794219488933794
    "G17" = 1.234295333333334
    "H17" = 3.702886
    "I17" = 0.01108709435540116
    "J17" = 0.01108709435540116
    "M17" = 65.761079
    "N17" = 197.283237
    "O17" = 0.1293522275572212
    "P17" = 0.1293522275572212
    "Q17" = 81.16859292466468
    "R17" = 730.517336321982
    "S17" = 0.001434140352008234
    "T17" = 0.001434140352008234
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

